# Add a new "segments" column before the existing data, and re-populate
# column A with a numeric 0-based index for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# RawActivations/PercActivations/totalActivation columns from B/C/D to C/D/E,
# and leaves column A (currently holding the segment names) untouched.
$ws.Columns.Item(2).Insert()

# The insert operation leaves stray formatting behind in the new column's
# data rows (copied from neighbouring cells), and no header formatting on
# row 1, so start from a clean slate across the whole inserted column.
$ws.Range("B1:B20").ClearFormats()

# New header for the inserted column, styled like the other header cells.
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Determine how many data rows there are (rows below the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Move the segment-name labels from column A into the new column B, and
# replace column A with a 0-based numeric index for each row.
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Text
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 1).Value = $r - 2
}
